$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix header row (correct "divison" typo to "division")
$ws.Range("E1").Value = "division"

# Remove the existing hyperlink on D2 so it can be re-added cleanly for every row
$ws.Hyperlinks.Delete() | Out-Null

# Student data for rows 2-11 (columns: sl_no, name, usn, email, division, gender)
$data = @(
  @(1,  "ABC", "01FS24BEC015", "01FS24BEC017@kletech.ac.in", "A", "female"),
  @(2,  "BCD", "01FS24BEC016", "01FS24BEC018@kletech.ac.in", "A", "female"),
  @(3,  "CDE", "01FS24BEC017", "01FS24BEC019@kletech.ac.in", "A", "female"),
  @(4,  "DEF", "01FS24BEC018", "01FS24BEC020@kletech.ac.in", "A", "male"),
  @(5,  "EFG", "01FS24BEC019", "01FS24BEC021@kletech.ac.in", "A", "female"),
  @(6,  "FGH", "01FS24BEC020", "01FS24BEC022@kletech.ac.in", "A", "female"),
  @(7,  "GHI", "01FS24BEC021", "01FS24BEC023@kletech.ac.in", "A", "female"),
  @(8,  "HIJ", "01FS24BEC022", "01FS24BEC024@kletech.ac.in", "A", "male"),
  @(9,  "IJK", "01FS24BEC023", "01FS24BEC025@kletech.ac.in", "A", "female"),
  @(10, "JKL", "01FS24BEC024", "01FS24BEC026@kletech.ac.in", "A", "female")
)

$row = 2
foreach ($rec in $data) {
    $ws.Cells.Item($row, 1).Value = $rec[0]
    $ws.Cells.Item($row, 2).Value = $rec[1]
    $ws.Cells.Item($row, 3).Value = $rec[2]

    $email = $rec[3]
    $dcell = $ws.Cells.Item($row, 4)
    $dcell.Value = $email
    $ws.Hyperlinks.Add($dcell, "mailto:$email") | Out-Null

    $ws.Cells.Item($row, 5).Value = $rec[4]
    $ws.Cells.Item($row, 6).Value = $rec[5]
    $row = $row + 1
}

# Adjust column B width to fit the longer header content
$ws.Columns.Item(2).ColumnWidth = 11.5

# Update the active selection
$ws.Range("I10").Select() | Out-Null
